$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "id_scenario" column (column A) is unnecessary and is being removed
# from the table. Every other column shifts one place to the left.

# Select column A first, mirroring the usual "right-click column header ->
# Delete" workflow in Excel.
[void]$ws.Columns.Item(1).Select()

# Work out how many columns currently hold data on row 1 so the shift
# below covers the whole used range regardless of sheet size.
$lastCol = $ws.Cells.Item(1, $ws.Columns.Count).End(-4159).Column

# Remember the explicit "best fit" widths of the columns that have one
# (everything except the id_region column and the trailing unit/number
# columns) before we start overwriting cell values, so we can re-apply
# them one position to the left afterwards.
$widths = @{}
for ($c = 3; $c -lt ($lastCol - 1); $c++) {
    $widths[$c] = $ws.Columns.Item($c).ColumnWidth
}

# Shift every column's row-1 value one position to the left (B->A, C->B, ...)
for ($c = 1; $c -lt $lastCol; $c++) {
    $src = $ws.Cells.Item(1, $c + 1)
    $dst = $ws.Cells.Item(1, $c)
    $dst.Value = $src.Value()
}

# Remove the now-duplicated trailing cell that used to hold the last column.
$ws.Cells.Item(1, $lastCol).Clear()

# Re-apply the remembered "best fit" widths to their new (shifted left)
# column position so the visual layout matches the old one.
foreach ($c in $widths.Keys) {
    $ws.Columns.Item($c - 1).ColumnWidth = $widths[$c]
}
